$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44243
$ws.Range("J2").Value = 1600
$ws.Range("K2").Value = 450
$ws.Range("L2").Value = 480
$ws.Range("M2").Value = 465
$ws.Range("O2").Value = 'Región del Maule'
$ws.Range("P2").Value = 465

# Row 3
$ws.Range("D3").Value = 44238
$ws.Range("J3").Value = 1250
$ws.Range("K3").Value = 430
$ws.Range("L3").Value = 450
$ws.Range("M3").Value = 440
$ws.Range("P3").Value = 440

# Row 4
$ws.Range("D4").Value = 44547
$ws.Range("I4").Value = '1a (cosecha)'
$ws.Range("J4").Value = 800
$ws.Range("K4").Value = 600
$ws.Range("L4").Value = 650
$ws.Range("M4").Value = 625
$ws.Range("O4").Value = 'Perú'
$ws.Range("P4").Value = 625

# Row 5
$ws.Range("D5").Value = 44547
$ws.Range("I5").Value = '2a nueva(o)'
$ws.Range("J5").Value = 300
$ws.Range("O5").Value = 'Perú'

# Row 6
$ws.Range("D6").Value = 44469
$ws.Range("I6").Value = '1a nueva(o)'
$ws.Range("J6").Value = 1200
$ws.Range("K6").Value = 600
$ws.Range("L6").Value = 650
$ws.Range("M6").Value = 625
$ws.Range("O6").Value = 'Perú'
$ws.Range("P6").Value = 625

# Row 7
$ws.Range("D7").Value = 44428
$ws.Range("I7").Value = '1a nueva(o)'
$ws.Range("J7").Value = 600
$ws.Range("K7").Value = 580
$ws.Range("L7").Value = 600
$ws.Range("M7").Value = 590
$ws.Range("P7").Value = 590

# Row 8
$ws.Range("D8").Value = 44399
$ws.Range("I8").Value = '1a (guarda)'
$ws.Range("J8").Value = 800
$ws.Range("K8").Value = 450
$ws.Range("L8").Value = 480
$ws.Range("M8").Value = 465
$ws.Range("O8").Value = 'Provincia de Melipilla'
$ws.Range("P8").Value = 465

# Row 9
$ws.Range("D9").Value = 44231
$ws.Range("J9").Value = 1300
$ws.Range("K9").Value = 450
$ws.Range("L9").Value = 480
$ws.Range("M9").Value = 465
$ws.Range("P9").Value = 465

# Row 11
$ws.Range("D11").Value = 44175
$ws.Range("J11").Value = 1200
$ws.Range("K11").Value = 1500
$ws.Range("L11").Value = 1550
$ws.Range("M11").Value = 1525
$ws.Range("P11").Value = 1525

# Row 12
$ws.Range("D12").Value = 44349
$ws.Range("H12").Value = 'Pachia'
$ws.Range("J12").Value = 1200
$ws.Range("K12").Value = 730
$ws.Range("L12").Value = 750
$ws.Range("M12").Value = 740
$ws.Range("P12").Value = 740

# Row 13
$ws.Range("D13").Value = 44453
$ws.Range("H13").Value = 'Camote'
$ws.Range("J13").Value = 800
$ws.Range("K13").Value = 630
$ws.Range("L13").Value = 650
$ws.Range("M13").Value = 640
$ws.Range("P13").Value = 640

# Row 14
$ws.Range("D14").Value = 44201
$ws.Range("I14").Value = '1a nueva(o)'
$ws.Range("J14").Value = 1360
$ws.Range("K14").Value = 730
$ws.Range("L14").Value = 750
$ws.Range("M14").Value = 740
$ws.Range("O14").Value = 'Perú'
$ws.Range("P14").Value = 740

# Row 15
$ws.Range("D15").Value = 44204
$ws.Range("I15").Value = '2a nueva(o)'
$ws.Range("J15").Value = 1600
$ws.Range("K15").Value = 500
$ws.Range("L15").Value = 550
$ws.Range("M15").Value = 525
$ws.Range("O15").Value = 'Región del Maule'
$ws.Range("P15").Value = 525

# Row 16
$ws.Range("D16").Value = 44168
$ws.Range("K16").Value = 1500
$ws.Range("L16").Value = 1700
$ws.Range("M16").Value = 1600
$ws.Range("P16").Value = 1600

# Row 17
$ws.Range("D17").Value = 44179
$ws.Range("J17").Value = 1000
$ws.Range("K17").Value = 1350
$ws.Range("L17").Value = 1400
$ws.Range("M17").Value = 1375
$ws.Range("P17").Value = 1375

# Row 18
$ws.Range("D18").Value = 44211
$ws.Range("J18").Value = 1600
$ws.Range("K18").Value = 500
$ws.Range("L18").Value = 550
$ws.Range("M18").Value = 525
$ws.Range("O18").Value = 'Región de O''Higgins'
$ws.Range("P18").Value = 525

# Row 19
$ws.Range("D19").Value = 44483
$ws.Range("J19").Value = 1300
$ws.Range("K19").Value = 550
$ws.Range("L19").Value = 580
$ws.Range("M19").Value = 565
$ws.Range("P19").Value = 565

# Row 20
$ws.Range("D20").Value = 44476
$ws.Range("I20").Value = '1a nueva(o)'
$ws.Range("K20").Value = 480
$ws.Range("L20").Value = 500
$ws.Range("M20").Value = 490
$ws.Range("O20").Value = 'Perú'
$ws.Range("P20").Value = 490

# Row 21
$ws.Range("D21").Value = 44301
$ws.Range("I21").Value = '2a nueva(o)'
$ws.Range("J21").Value = 1200
$ws.Range("K21").Value = 400
$ws.Range("L21").Value = 430
$ws.Range("M21").Value = 415
$ws.Range("O21").Value = 'Provincia de Melipilla'
$ws.Range("P21").Value = 415

# Row 22
$ws.Range("D22").Value = 44322
$ws.Range("I22").Value = '1a (cosecha)'
$ws.Range("K22").Value = 350
$ws.Range("L22").Value = 400
$ws.Range("M22").Value = 375
$ws.Range("O22").Value = 'Región del Maule'
$ws.Range("P22").Value = 375

# Row 23
$ws.Range("D23").Value = 44490
$ws.Range("J23").Value = 1200
$ws.Range("K23").Value = 450
$ws.Range("L23").Value = 480
$ws.Range("M23").Value = 465
$ws.Range("P23").Value = 465
